# Update "Top Gainers" rows 52-76: the top-of-range row (old MIDWESTLTD @52)
# drops off the list, every following row shifts up by one, CRAMC's Weekly
# value is corrected to 7.7681, and a new row (DPSCLTD) is appended at the
# bottom (row 76).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top Gainers")

$rows = @(
    @{ Row = 52; Stock = "BHEL";       Latest = 1.8175; Weekly = 8.1883;  Monthly = 4.7809 },
    @{ Row = 53; Stock = "THELEELA";   Latest = 1.7961; Weekly = 1.7961;  Monthly = 5.6659 },
    @{ Row = 54; Stock = "DPABHUSHAN"; Latest = 1.7774; Weekly = 0.4476;  Monthly = -2.6106 },
    @{ Row = 55; Stock = "EXICOM";     Latest = 1.7692; Weekly = 1.9152;  Monthly = -1.5316 },
    @{ Row = 56; Stock = "PNCINFRA";   Latest = 1.769;  Weekly = 1.679;   Monthly = -2.1599 },
    @{ Row = 57; Stock = "FISCHER";    Latest = 1.7685; Weekly = 12.6908; Monthly = 5.7568 },
    @{ Row = 58; Stock = "CRAMC";      Latest = 1.6891; Weekly = 7.7681;  Monthly = "N/A" },
    @{ Row = 59; Stock = "CENTRUM";    Latest = 1.6667; Weekly = 0.2989;  Monthly = -0.3564 },
    @{ Row = 60; Stock = "MARINE";     Latest = 1.6655; Weekly = -1.6528; Monthly = 10.3569 },
    @{ Row = 61; Stock = "OIL";        Latest = 1.6647; Weekly = 1.907;   Monthly = 3.3108 },
    @{ Row = 62; Stock = "IFCI";       Latest = 1.6296; Weekly = 4.1854;  Monthly = 7.308 },
    @{ Row = 63; Stock = "SEQUENT";    Latest = 1.6269; Weekly = 6.9633;  Monthly = 16.018 },
    @{ Row = 64; Stock = "NETWEB";     Latest = 1.6233; Weekly = 7.2852;  Monthly = 9.3025 },
    @{ Row = 65; Stock = "JKTYRE";     Latest = 1.6032; Weekly = 4.6122;  Monthly = 20.4949 },
    @{ Row = 66; Stock = "IIFLCAPS";   Latest = 1.6019; Weekly = -4.4364; Monthly = 22.387 },
    @{ Row = 67; Stock = "SULA";       Latest = 1.5786; Weekly = -1.2202; Monthly = -3.6216 },
    @{ Row = 68; Stock = "SCHNEIDER";  Latest = 1.5706; Weekly = 3.1055;  Monthly = 5.5996 },
    @{ Row = 69; Stock = "NRBBEARING"; Latest = 1.5385; Weekly = 1.3293;  Monthly = -4.6344 },
    @{ Row = 70; Stock = "MOTILALOFS"; Latest = 1.533;  Weekly = -1.2498; Monthly = 14.5087 },
    @{ Row = 71; Stock = "IMPAL";      Latest = 1.5245; Weekly = 2.1;     Monthly = 4.7574 },
    @{ Row = 72; Stock = "SAURASHCEM"; Latest = 1.5225; Weekly = 0.9536;  Monthly = -1.6285 },
    @{ Row = 73; Stock = "NAZARA";     Latest = 1.5174; Weekly = 0.9500999999999999; Monthly = 6.4218 },
    @{ Row = 74; Stock = "NEULANDLAB"; Latest = 1.454;  Weekly = -2.8468; Monthly = 7.0181 },
    @{ Row = 75; Stock = "ASTRAMICRO"; Latest = 1.4463; Weekly = -2.4533; Monthly = 7.2924 },
    @{ Row = 76; Stock = "DPSCLTD";    Latest = 1.4456; Weekly = 1.0161;  Monthly = -1.3234 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Stock
    $ws.Cells.Item($r.Row, 3).Value = $r.Latest
    $ws.Cells.Item($r.Row, 4).Value = $r.Weekly
    $ws.Cells.Item($r.Row, 5).Value = $r.Monthly
}
